# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracking sheet and
# moves the "latest row" date style down to the newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (34) used the distinct "date only" number format
# reserved for the most recent entry; now that row 35 is the newest entry,
# row 34 reverts to the regular "date + time" format shared by the rest of
# the table.
$ws.Range("A34").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 35.
$ws.Range("A35").Value = 45775
$ws.Range("B35").Value = 139
$ws.Range("C35").Value = 146
$ws.Range("D35").Value = 143

# The newest row now carries the distinct "date only" number format that
# previously marked row 34 as the latest entry.
$ws.Range("A35").NumberFormat = "YYYY-MM-DD"
